$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# --- Sheet 1 ---
$ws1.Range("A2").Value = "DTaP "
$ws1.Range("A3").Value = "DTaP "
$ws1.Range("A4").Value = "DTaP "
$ws1.Range("A5").Value = "DTaP "
$ws1.Range("A6").Value = "DTaP-IPV "
$ws1.Range("A7").Value = "DTaP-IPV "
$ws1.Range("A8").Value = "DTaP-Hep B-IPV "
$ws1.Range("A9").Value = "DTaP-IP-HI "
$ws1.Range("A10").Value = "e-IPV "
$ws1.Range("A11").Value = "Hepatitis B-Hib "
$ws1.Range("A12").Value = "Hepatitis A Pediatric "
$ws1.Range("A13").Value = "Hepatitis A Pediatric "
$ws1.Range("A14").Value = "Hepatitis A Pediatric "
$ws1.Range("A15").Value = "Hepatitis A-Hepatitis B 18 only "
$ws1.Range("A16").Value = "Hepatitis A-Hepatitis B 18 only "
$ws1.Range("A17").Value = "Hepatitis B  Pediatric/Adolescent"
$ws1.Range("A18").Value = "Hepatitis B  Pediatric/Adolescent"
$ws1.Range("A19").Value = "Hepatitis B  Pediatric/Adolescent"
$ws1.Range("B19").Value = "Recombivax HB"
$ws1.Range("A20").Value = "Hib "
$ws1.Range("A21").Value = "Hib "
$ws1.Range("A22").Value = "Hib "
$ws1.Range("A23").Value = "HPV - Quadrivalent Human Papillomavirus Types 6, 11, 16 and 18 Recombinant "
$ws1.Range("A24").Value = "HPV -Bivalent Human Papillomavirus Types 16 and 18 "
$ws1.Range("A25").Value = "HPV -Bivalent Human Papillomavirus Types 16 and 18 "
$ws1.Range("A26").Value = "Measles, Mumps, Rubella and Varicella (MMR-V) "
$ws1.Range("A27").Value = "Meningococcal Conjugate (Groups A, C, Y and W-135) "
$ws1.Range("A28").Value = "Meningococcal Conjugate (Groups A, C, Y and W-135) "
$ws1.Range("A29").Value = "Measles, Mumps and Rubella (MMR) "
$ws1.Range("A30").Value = "Pneumococcal 13-valent  (Pediatric)"
$ws1.Range("A32").Value = "Rotavirus, Live, Oral, Pentavalent "
$ws1.Range("A33").Value = "Rotavirus, Live, Oral, Oral "
$ws1.Range("A34").Value = "Tetanus  Diphtheria Toxoids "
$ws1.Range("D34").Value = "10 pack - 1 dose syringes No Needle"
$ws1.Range("A35").Value = "Tetanus  Diphtheria Toxoids "
$ws1.Range("A36").Value = "Tetanus  Diphtheria Toxoids "
$ws1.Range("A37").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws1.Range("A38").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws1.Range("D38").Value = "10 pack - 1 dose TL syringes, No Needle"
$ws1.Range("A39").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws1.Range("A40").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws1.Range("A41").Value = "Varicella "

# --- Sheet 2 ---
$ws2.Range("A2").Value = "Hepatitis A Adult "
$ws2.Range("A3").Value = "Hepatitis A Adult "
$ws2.Range("A4").Value = "Hepatitis A-Hepatitis B Adult "
$ws2.Range("A5").Value = "Hepatitis A-Hepatitis B Adult "
$ws2.Range("A6").Value = "Hepatitis B-Adult "
$ws2.Range("A7").Value = "Hepatitis B-Adult "
$ws2.Range("A8").Value = "HPV -Quadrivalent Human Papillomavirus Types 6, 11, 16 and 18 Recombinant Adult "
$ws2.Range("A9").Value = "HPV-Human Papillomavirus Bivalent Types 16 and 18 "
$ws2.Range("A10").Value = "HPV-Human Papillomavirus Bivalent Types 16 and 18 "
$ws2.Range("A11").Value = "Measles, Mumps,  Rubella-Adult "
$ws2.Range("A14").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws2.Range("A15").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws2.Range("A16").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws2.Range("A17").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws2.Range("A18").Value = "Varicella-Adult "
$ws2.Range("A21").Value = "Tetanus and Diphtheria Toxoids "
$ws2.Range("A22").Value = "Meningococcal Conjugate (Groups A, C, H and W-135) "

# --- Sheet 3 ---
$ws3.Range("A2").Value = "Influenza  (Age 6 months and older)"
$ws3.Range("A3").Value = "Influenza  (Age 6-35 months)"
$ws3.Range("B3").Value = "Fluzone Pediatric dose No Preservative"
$ws3.Range("A4").Value = "Influenza  (Age 36 months and older)"
$ws3.Range("B4").Value = "Fluzone No-Preservative"
$ws3.Range("A5").Value = "Influenza  (Age 36 months and older)"
$ws3.Range("B5").Value = "Fluzone No-Preservative"
$ws3.Range("A6").Value = "Influenza  (Age 36 months and older)"
$ws3.Range("B6").Value = "Fluarix Preservative Free"
$ws3.Range("D6").Value = "10 pack- 1 dose TipLok syringe"
$ws3.Range("A7").Value = "Influenza  (Age 4 years and older)"
$ws3.Range("A8").Value = "Influenza  (Age 4 years and older)"
$ws3.Range("B8").Value = "Fluvirin Preservative Free"
$ws3.Range("A9").Value = "Influenza  Live, Intranasal (Age 2-49 years)"
$ws3.Range("B9").Value = "FluMist No Preservative"
$ws3.Range("A10").Value = "Influenza  (Age 36 months and older)"
$ws3.Range("B10").Value = "Afluria No Preservative"
$ws3.Range("D10").Value = "10 pack-1 dose syringe"
$ws3.Range("H10").Value = "Merck (CSL product)"

# --- Sheet 4 ---
$ws4.Range("A2").Value = "Influenza "
$ws4.Range("A3").Value = "Influenza "
$ws4.Range("B3").Value = "Fluzone No Preservative"
$ws4.Range("A4").Value = "Influenza "
$ws4.Range("B4").Value = "Fluzone No Preservative"
$ws4.Range("A5").Value = "Influenza "
$ws4.Range("A6").Value = "Influenza "
$ws4.Range("A7").Value = "Influenza "
$ws4.Range("B7").Value = "Fluvirin Preservative Free"
$ws4.Range("A8").Value = "Influenza "
$ws4.Range("B8").Value = "Fluarix Preservative Free"
$ws4.Range("A9").Value = "Influenza "
$ws4.Range("A10").Value = "Influenza "
$ws4.Range("B10").Value = "FluMist No Preservative"

Write-Host "Done applying footnote/line-break cleanup."